# Apply data corrections to p6_analysis.xlsx
# "sunday data rain is falling" - updating FeltMotion coded responses on Sheet1
# and the corresponding aggregated probabilities on FeltMotion_Prob.

$wb = $excel.ActiveWorkbook

# --- Sheet1: raw trial-level coding corrections ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("G2").Value  = 0

$ws1.Range("G10").Value = 1

$ws1.Range("F47").Value = 1
$ws1.Range("G47").Value = 1
$ws1.Range("I47").Value = 0

$ws1.Range("G51").Value = 1

$ws1.Range("F63").Value = 1
$ws1.Range("I63").Value = 0

$ws1.Range("F77").Value = 0
$ws1.Range("G77").Value = 1
$ws1.Range("I77").Value = 1

$ws1.Range("G91").Value = 1

$ws1.Range("G119").Value = 1

$ws1.Range("G132").Value = 1

$ws1.Range("F145").Value = 1
$ws1.Range("I145").Value = 0

$ws1.Range("G148").Value = 1

$ws1.Range("G163").Value = 1

$ws1.Range("G177").Value = 1

# --- FeltMotion_Prob: recomputed aggregate mean/sem values ---
$ws2 = $wb.Worksheets.Item("FeltMotion_Prob")

$ws2.Range("D5").Value = 0.3333333333333333
$ws2.Range("E5").Value = 0.1666666666666667

$ws2.Range("D14").Value = 0.7142857142857143
$ws2.Range("E14").Value = 0.1844277783908294
